$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 41669452
$ws.Range("I18").Value = 50001344
$ws.Range("J18").Value = 10000
$ws.Range("K18").Value = 50001344
$ws.Range("L18").Value = 10000
$ws.Range("M18").Value = -50001060
$ws.Range("N18").Value = -10568
$ws.Range("H33").Value = 653.8148
$ws.Range("I33").Value = 707.04346
$ws.Range("K33").Value = 707.04346
$ws.Range("M33").Value = -478.04346
$ws.Range("H40").Value = 4799.9287
$ws.Range("I40").Value = 4799.9287
$ws.Range("K40").Value = 4799.9287
$ws.Range("M40").Value = -4624.9287
$ws.Range("H43").Value = 1444.6
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H47").Value = 47750
$ws.Range("J47").Value = 48666.668
$ws.Range("L47").Value = 48666.668
$ws.Range("N47").Value = -50610.668
$ws.Range("H64").Value = 4136.9
$ws.Range("J64").Value = 4944.25
$ws.Range("L64").Value = 4944.25
$ws.Range("N64").Value = -5440.25
$ws.Range("H67").Value = 4136.9
$ws.Range("J67").Value = 4944.25
$ws.Range("L67").Value = 4944.25
$ws.Range("N67").Value = -6660.25
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H88").Value = 83362830
$ws.Range("I88").Value = 166669660
$ws.Range("J88").Value = 55999.5
$ws.Range("K88").Value = 166669660
$ws.Range("L88").Value = 55999.5
$ws.Range("M88").Value = -166669254
$ws.Range("N88").Value = -56811.5
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H91").Value = 83362830
$ws.Range("I91").Value = 166669660
$ws.Range("J91").Value = 55999.5
$ws.Range("K91").Value = 166669660
$ws.Range("L91").Value = 55999.5
$ws.Range("M91").Value = -166668256
$ws.Range("N91").Value = -58807.5
$ws.Range("H93").Value = 23866
$ws.Range("J93").Value = 23866
$ws.Range("L93").Value = 23866
$ws.Range("N93").Value = -28858
$ws.Range("H100").Value = 1789.3
$ws.Range("I100").Value = 1387.8572
$ws.Range("K100").Value = 1387.8572
$ws.Range("M100").Value = -846.8571999999999
$ws.Range("H107").Value = 20011128
$ws.Range("I107").Value = 21750140
$ws.Range("J107").Value = 12499.5
$ws.Range("K107").Value = 21750140
$ws.Range("L107").Value = 12499.5
$ws.Range("M107").Value = -21748220
$ws.Range("N107").Value = -16339.5
$ws.Range("H112").Value = 1118.6538
$ws.Range("J112").Value = 1155.625
$ws.Range("L112").Value = 3466.875
$ws.Range("N112").Value = -5682.875
$ws.Range("H132").Value = 6949.4873
$ws.Range("I132").Value = 3232.8147
$ws.Range("K132").Value = 9698.444100000001
$ws.Range("M132").Value = -7168.444100000001
$ws.Range("H137").Value = 1290.2632
$ws.Range("I137").Value = 1081.875
$ws.Range("K137").Value = 3245.625
$ws.Range("M137").Value = -695.625

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6595.1377
$ws.Range("I74").Value = 6779.4614
$ws.Range("J74").Value = 4997.6665
$ws.Range("K74").Value = 6779.4614
$ws.Range("L74").Value = 4997.6665
$ws.Range("M74").Value = -5905.4614
$ws.Range("N74").Value = -6745.6665
$ws.Range("H77").Value = 6595.1377
$ws.Range("I77").Value = 6779.4614
$ws.Range("J77").Value = 4997.6665
$ws.Range("K77").Value = 33897.307
$ws.Range("L77").Value = 24988.3325
$ws.Range("M77").Value = -29529.307
$ws.Range("N77").Value = -33724.3325
$ws.Range("H110").Value = 2879.6086
$ws.Range("I110").Value = 2878.7273
$ws.Range("J110").Value = 2899
$ws.Range("K110").Value = 2878.7273
$ws.Range("L110").Value = 2899
$ws.Range("M110").Value = -833.7273
$ws.Range("N110").Value = -6989

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 467.27274
$ws.Range("I22").Value = 262.625
$ws.Range("J22").Value = 1013
$ws.Range("K22").Value = 262.625
$ws.Range("L22").Value = 1013
$ws.Range("M22").Value = -89.625
$ws.Range("N22").Value = -1359
$ws.Range("H99").Value = 3297.1538
$ws.Range("I99").Value = 3151.7273
$ws.Range("K99").Value = 3151.7273
$ws.Range("M99").Value = -1653.7273
$ws.Range("H134").Value = 7111.76
$ws.Range("I134").Value = 7889.737
$ws.Range("K134").Value = 23669.211
$ws.Range("M134").Value = -21134.211
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("N141").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3320.75
$ws.Range("I16").Value = 3827.6667
$ws.Range("K16").Value = 3827.6667
$ws.Range("M16").Value = -3540.6667
$ws.Range("H22").Value = 1809.7778
$ws.Range("I22").Value = 1683.8
$ws.Range("K22").Value = 1683.8
$ws.Range("M22").Value = -1333.8
$ws.Range("H113").Value = 3320.75
$ws.Range("I113").Value = 3827.6667
$ws.Range("K113").Value = 3827.6667
$ws.Range("M113").Value = -1657.6667
$ws.Range("H134").Value = 10846.849
$ws.Range("I134").Value = 11701.643
$ws.Range("J134").Value = 6060
$ws.Range("K134").Value = 35104.929
$ws.Range("L134").Value = 18180
$ws.Range("M134").Value = -32569.929
$ws.Range("N134").Value = -23250

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 14706563
$ws.Range("I92").Value = 83333950
$ws.Range("K92").Value = 250001850
$ws.Range("M92").Value = -250000602

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 27540.2
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H70").Value = 3765.7693
$ws.Range("I70").Value = 3138
$ws.Range("J70").Value = 4621.8184
$ws.Range("K70").Value = 3138
$ws.Range("L70").Value = 4621.8184
$ws.Range("M70").Value = -2868
$ws.Range("N70").Value = -5161.8184
$ws.Range("H73").Value = 3765.7693
$ws.Range("I73").Value = 3138
$ws.Range("J73").Value = 4621.8184
$ws.Range("K73").Value = 3138
$ws.Range("L73").Value = 4621.8184
$ws.Range("M73").Value = -2202
$ws.Range("N73").Value = -6493.8184

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11662.721
$ws.Range("I40").Value = 11686.2705
$ws.Range("K40").Value = 11686.2705
$ws.Range("M40").Value = -11550.2705
$ws.Range("H68").Value = 12822450
$ws.Range("I68").Value = 12822450
$ws.Range("K68").Value = 12822450
$ws.Range("M68").Value = -12821701
$ws.Range("H71").Value = 12822450
$ws.Range("I71").Value = 12822450
$ws.Range("K71").Value = 64112250
$ws.Range("M71").Value = -64108506

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 764.46155
$ws.Range("I100").Value = 668.625
$ws.Range("K100").Value = 1337.25
$ws.Range("M100").Value = -796.25
$ws.Range("H132").Value = 1521.2858
$ws.Range("I132").Value = 1279
$ws.Range("K132").Value = 3837
$ws.Range("M132").Value = -1307
